$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 2: insert a new "Age" (continuous stats) row above the existing
# "Sex" row, and relabel the "Sex" row's m/w to male/female.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")

# Insert a fresh blank row at row 2 (everything from old row 2 onward shifts
# down by one); strip the formatting it inherited from the header row.
$ws2.Rows.Item(2).Insert()
$ws2.Rows.Item(2).ClearFormats()

# New row 2: continuous "Age" statistics.
$ws2.Range("A2").Value = "Age"
$ws2.Range("B2").Value = "mean(SD) = 58.8 (13.1)`nmedian(IQR) = 57 (51 - 70)`nrange = 20 - 89`ncomplete: n = 91"

# Row 3 (previously row 2): "Sex" row, relabel m/w -> male/female.
$ws2.Range("A3").Value = "Sex"
$ws2.Range("B3").Value = "male: 61.5% (56)`nfemale: 38.5% (35)`ncomplete: n = 91"

# Restore default (auto) row heights - the multi-line values above would
# otherwise leave an explicit customHeight on these two rows.
$ws2.Rows.Item(2).AutoFit()
$ws2.Rows.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# Table 3: same change, across the 4 severity columns (B..E) plus the
# significance column (F).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).ClearFormats()

# New row 2: continuous "Age" statistics per column.
$ws3.Range("A2").Value = "Age"
$ws3.Range("B2").Value = "mean(SD) = 48.6 (12.5)`nmedian(IQR) = 50 (41.5 - 56.5)`nrange = 20 - 72`ncomplete: n = 19"
$ws3.Range("C2").Value = "mean(SD) = 60.6 (12.7)`nmedian(IQR) = 61 (53 - 71.5)`nrange = 38 - 81`ncomplete: n = 23"
$ws3.Range("D2").Value = "mean(SD) = 64.7 (13.5)`nmedian(IQR) = 64 (56.5 - 77)`nrange = 37 - 89`ncomplete: n = 23"
$ws3.Range("E2").Value = "mean(SD) = 59.4 (9.38)`nmedian(IQR) = 56.5 (53.2 - 65.8)`nrange = 45 - 80`ncomplete: n = 26"
$ws3.Range("F2").Value = "p = 0.0019"

# Row 3 (previously row 2): "Sex" row, relabel m/w -> male/female.
$ws3.Range("A3").Value = "Sex"
$ws3.Range("B3").Value = "male: 31.6% (6)`nfemale: 68.4% (13)`ncomplete: n = 19"
$ws3.Range("C3").Value = "male: 47.8% (11)`nfemale: 52.2% (12)`ncomplete: n = 23"
$ws3.Range("D3").Value = "male: 87% (20)`nfemale: 13% (3)`ncomplete: n = 23"
$ws3.Range("E3").Value = "male: 73.1% (19)`nfemale: 26.9% (7)`ncomplete: n = 26"
$ws3.Range("F3").Value = "p = 0.00079"

# Restore default (auto) row heights - the multi-line values above would
# otherwise leave an explicit customHeight on these two rows.
$ws3.Rows.Item(2).AutoFit()
$ws3.Rows.Item(3).AutoFit()
